$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.982.40"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "3.522.99"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.618"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("D8").Value = "3.519.83"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.200"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("E12").Value = "  -2.99%  "
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000279"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "4.089.46"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "622.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.94%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.524.11"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "70.059.79"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.886"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -11.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  -2.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("E32").Value = "  -4.10%  "
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "568.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.49%  "
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "56.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0453"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.65%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.141"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.13%  "
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("D44").Value = "3.341.43"
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "33.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.17%  "
$ws.Range("B47").Value = "PEPE"
$ws.Range("C47").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D47").Value = "0.0₃0706"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("E49").Value = "  -3.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.61%  "
$ws.Range("E51").Value = "  +4.52%  "
